$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Copy the header cell's formatting (bold, border, centered) onto the new
# "Serviced by " header cell before we give it its own value.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("O1").Value = "Serviced by "

# Drop the trailing space that used to be on the "Correction " header.
$ws.Range("N1").Value = "Correction"

# The "Correction" column used to be left completely blank on every data
# row; now it is filled in with the same "nan" placeholder text used by
# the other empty columns.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

Write-Output "Added 'Serviced by ' column to Card23"
